# Apply the "Add files via upload" edit to the Test-Cases sheet:
#  - I2 (Approved/Rejected column) flips from "Approved" to "Rejected"
#  - J2 (ReasonToReject column) gains the value "unmatched"
#  - The active selection moves from I4 to K2

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test-Cases")

$ws.Range("I2").Value = "Rejected"
$ws.Range("J2").Value = "unmatched"

$ws.Range("K2").Select()
